$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5341.6665
$ws.Range("J18").Value = 14746.75
$ws.Range("L18").Value = 14746.75
$ws.Range("N18").Value = -15314.75
$ws.Range("H29").Value = 66.333336
$ws.Range("I29").Value = 66.333336
$ws.Range("K29").Value = 199.000008
$ws.Range("M29").Value = 81.99999199999999
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H51").Value = 8446.32
$ws.Range("J51").Value = 8399.044
$ws.Range("L51").Value = 8399.044
$ws.Range("N51").Value = -9367.044
$ws.Range("H55").Value = 251.5
$ws.Range("J55").Value = 662.8
$ws.Range("L55").Value = 662.8
$ws.Range("N55").Value = -1090.8
$ws.Range("H62").Value = 8780148
$ws.Range("I62").Value = 11119130
$ws.Range("J62").Value = 8967.75
$ws.Range("K62").Value = 11119130
$ws.Range("L62").Value = 8967.75
$ws.Range("M62").Value = -11118506
$ws.Range("N62").Value = -10215.75
$ws.Range("H65").Value = 8780148
$ws.Range("I65").Value = 11119130
$ws.Range("J65").Value = 8967.75
$ws.Range("K65").Value = 55595650
$ws.Range("L65").Value = 44838.75
$ws.Range("M65").Value = -55592530
$ws.Range("N65").Value = -51078.75
$ws.Range("H74").Value = 4957
$ws.Range("I74").Value = 4957
$ws.Range("K74").Value = 4957
$ws.Range("M74").Value = -4021
$ws.Range("H77").Value = 4957
$ws.Range("I77").Value = 4957
$ws.Range("K77").Value = 24785
$ws.Range("M77").Value = -20105
$ws.Range("H98").Value = 715791.1
$ws.Range("I98").Value = 1589.6666
$ws.Range("K98").Value = 1589.6666
$ws.Range("M98").Value = -91.66660000000002
$ws.Range("H117").Value = 69999
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 69999
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 69999
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -79177
$ws.Range("H122").Value = 715791.1
$ws.Range("I122").Value = 1589.6666
$ws.Range("K122").Value = 4768.9998
$ws.Range("M122").Value = -2318.9998
$ws.Range("H132").Value = 3164.5386
$ws.Range("I132").Value = 2584.7097
$ws.Range("K132").Value = 7754.1291
$ws.Range("M132").Value = -5224.1291
$ws.Range("H138").Value = 2840.4695
$ws.Range("I138").Value = 1589.1
$ws.Range("J138").Value = 3161.3333
$ws.Range("K138").Value = 4767.299999999999
$ws.Range("L138").Value = 9483.999899999999
$ws.Range("M138").Value = 372.7000000000007
$ws.Range("N138").Value = -19763.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2461.7288
$ws.Range("I32").Value = 2461.7288
$ws.Range("K32").Value = 2461.7288
$ws.Range("M32").Value = -2174.7288
$ws.Range("H74").Value = 78406.16
$ws.Range("I74").Value = 99741.7
$ws.Range("J74").Value = 7287.6665
$ws.Range("K74").Value = 99741.7
$ws.Range("L74").Value = 7287.6665
$ws.Range("M74").Value = -98867.7
$ws.Range("N74").Value = -9035.666499999999
$ws.Range("H77").Value = 78406.16
$ws.Range("I77").Value = 99741.7
$ws.Range("J77").Value = 7287.6665
$ws.Range("K77").Value = 498708.5
$ws.Range("L77").Value = 36438.3325
$ws.Range("M77").Value = -494340.5
$ws.Range("N77").Value = -45174.3325
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H117").Value = 143415.67
$ws.Range("J117").Value = 143415.67
$ws.Range("L117").Value = 143415.67
$ws.Range("N117").Value = -152593.67

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 39148.168
$ws.Range("J81").Value = 36977.8
$ws.Range("L81").Value = 36977.8
$ws.Range("N81").Value = -39099.8
$ws.Range("H84").Value = 39148.168
$ws.Range("J84").Value = 36977.8
$ws.Range("L84").Value = 110933.4
$ws.Range("N84").Value = -121541.4
$ws.Range("H105").Value = 1983.5758
$ws.Range("I105").Value = 1930.6333
$ws.Range("J105").Value = 2513
$ws.Range("K105").Value = 1930.6333
$ws.Range("L105").Value = 2513
$ws.Range("M105").Value = -183.6333
$ws.Range("N105").Value = -6007
$ws.Range("H107").Value = 825.7143
$ws.Range("I107").Value = 825.7143
$ws.Range("K107").Value = 825.7143
$ws.Range("M107").Value = 1094.2857
$ws.Range("H138").Value = 60067.043
$ws.Range("J138").Value = 60067.043
$ws.Range("L138").Value = 60067.043
$ws.Range("N138").Value = -70347.04300000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2656.7144
$ws.Range("I58").Value = 2070.4
$ws.Range("J58").Value = 4122.5
$ws.Range("K58").Value = 2070.4
$ws.Range("L58").Value = 4122.5
$ws.Range("M58").Value = -1867.4
$ws.Range("N58").Value = -4528.5
$ws.Range("H99").Value = 366872.66
$ws.Range("I99").Value = 1253669.9
$ws.Range("J99").Value = 12153.75
$ws.Range("K99").Value = 1253669.9
$ws.Range("L99").Value = 12153.75
$ws.Range("M99").Value = -1252171.9
$ws.Range("N99").Value = -15149.75
$ws.Range("H107").Value = 3269.6785
$ws.Range("J107").Value = 5049.758
$ws.Range("L107").Value = 5049.758
$ws.Range("N107").Value = -8889.758
$ws.Range("H122").Value = 2127.5833
$ws.Range("I122").Value = 1948.1111
$ws.Range("K122").Value = 5844.3333
$ws.Range("M122").Value = -3394.3333
$ws.Range("H126").Value = 366872.66
$ws.Range("I126").Value = 1253669.9
$ws.Range("J126").Value = 12153.75
$ws.Range("K126").Value = 3761009.7
$ws.Range("L126").Value = 36461.25
$ws.Range("M126").Value = -3758539.7
$ws.Range("N126").Value = -41401.25
$ws.Range("H136").Value = 2656.7144
$ws.Range("I136").Value = 2070.4
$ws.Range("J136").Value = 4122.5
$ws.Range("K136").Value = 6211.200000000001
$ws.Range("L136").Value = 12367.5
$ws.Range("M136").Value = -3661.200000000001
$ws.Range("N136").Value = -17467.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1115.25
$ws.Range("I107").Value = 397.2857
$ws.Range("K107").Value = 1191.8571
$ws.Range("M107").Value = 728.1428999999998
$ws.Range("H127").Value = 1805.6
$ws.Range("J127").Value = 1633.25
$ws.Range("L127").Value = 4899.75
$ws.Range("N127").Value = -14819.75
$ws.Range("H131").Value = 10418305
$ws.Range("I131").Value = 83334550
$ws.Range("K131").Value = 250003650
$ws.Range("M131").Value = -249998610
$ws.Range("H136").Value = 2943.7144
$ws.Range("I136").Value = 2943.7144
$ws.Range("K136").Value = 8831.143199999999
$ws.Range("M136").Value = -3731.143199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 125003560
$ws.Range("I70").Value = 4067.8572
$ws.Range("K70").Value = 4067.8572
$ws.Range("M70").Value = -3797.8572
$ws.Range("H73").Value = 125003560
$ws.Range("I73").Value = 4067.8572
$ws.Range("K73").Value = 4067.8572
$ws.Range("M73").Value = -3131.8572
$ws.Range("H97").Value = 1497.0416
$ws.Range("I97").Value = 1325.238
$ws.Range("K97").Value = 1325.238
$ws.Range("M97").Value = -829.2380000000001
$ws.Range("H98").Value = 17534
$ws.Range("J98").Value = 17534
$ws.Range("L98").Value = 17534
$ws.Range("N98").Value = -23524
$ws.Range("H102").Value = 24197.348
$ws.Range("J102").Value = 103514.9
$ws.Range("L102").Value = 103514.9
$ws.Range("N102").Value = -106758.9
$ws.Range("H113").Value = 6341.857
$ws.Range("I113").Value = 10048.75
$ws.Range("K113").Value = 10048.75
$ws.Range("M113").Value = -7878.75
$ws.Range("H126").Value = 5127
$ws.Range("I126").Value = 5164.36
$ws.Range("K126").Value = 15493.08
$ws.Range("M126").Value = -13023.08
$ws.Range("H132").Value = 25064.977
$ws.Range("I132").Value = 31005.234
$ws.Range("J132").Value = 2624
$ws.Range("K132").Value = 93015.702
$ws.Range("L132").Value = 7872
$ws.Range("M132").Value = -90485.702
$ws.Range("N132").Value = -12932
$ws.Range("H136").Value = 28414.084
$ws.Range("J136").Value = 28414.084
$ws.Range("L136").Value = 85242.25199999999
$ws.Range("N136").Value = -90342.25199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5572.9
$ws.Range("I7").Value = 4759.6
$ws.Range("J7").Value = 6386.2
$ws.Range("K7").Value = 4759.6
$ws.Range("L7").Value = 6386.2
$ws.Range("M7").Value = -4647.6
$ws.Range("N7").Value = -6610.2
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H101").Value = 13666.667
$ws.Range("J101").Value = 13666.667
$ws.Range("L101").Value = 13666.667
$ws.Range("N101").Value = -20156.667
$ws.Range("H122").Value = 2748
$ws.Range("I122").Value = 2297.6
$ws.Range("K122").Value = 6892.799999999999
$ws.Range("M122").Value = -4442.799999999999
$ws.Range("H126").Value = 5572.9
$ws.Range("I126").Value = 4759.6
$ws.Range("J126").Value = 6386.2
$ws.Range("K126").Value = 14278.8
$ws.Range("L126").Value = 19158.6
$ws.Range("M126").Value = -11808.8
$ws.Range("N126").Value = -24098.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4898.8
$ws.Range("I81").Value = 4898.8
$ws.Range("K81").Value = 9797.6
$ws.Range("M81").Value = -8736.6
$ws.Range("H84").Value = 4898.8
$ws.Range("I84").Value = 4898.8
$ws.Range("K84").Value = 48988
$ws.Range("M84").Value = -43684
